$d = $word.ActiveDocument

$xmlBody = @'
<w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t>Java Programming Project – Asteroids</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Configuration:</w:t></w:r></w:p><w:p><w:r><w:t>Introduce a default configuration file (inside a jar).  Introduce a user configuration file (filesystem).  Add a class to load and manage all that stuff.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Acceleration/Deceleration:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">First attempt at motion just moves ship in given direction.  Update to include speed increments </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>Logging????</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:r><w:t xml:space="preserve">Do we introduce Log4j etc.?  </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Useful?</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>KeyReleases</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>??</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:r><w:t xml:space="preserve">Deal with </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>keyReleases</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">? Set </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>boolean</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> flag on key press so that, for example, we don't attempt to play a sound a second time until the key releases?  Can we detect if "the sound has ended"?? Address</w:t></w:r><w:r><w:t xml:space="preserve"> these issues</w:t></w:r></w:p><w:p/><w:p/><w:p/>
'@

$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $xmlBody + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$d.Content.InsertXML($xmlFrag)

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
